$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: Hydrogen demand updated; Non-metallic minerals value cleared
$ws.Range("B3").Value = 51248973.67504942
$ws.Range("D3").Value = ""

# Row 4: Methanol / Chemicals value corrected
$ws.Range("C4").Value = 8683.012753149565

# Row 5: Ammonia / Chemicals value corrected
$ws.Range("C5").Value = 13028.31263971192

# Row 7: renamed from "Other" to "Biogas", value corrected
$ws.Range("A7").Value = "Biogas"
$ws.Range("D7").Value = 24607.06425872971

# New row 8: "Other" category (previously row 7), with its own value
$ws.Range("A8").Value = "Other"
$ws.Range("A7").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("D8").Value = 16522.67512365857

$excel.CutCopyMode = 0
